$wb = $excel.ActiveWorkbook

# --- Add the new "25/04/2019" workout-log entries -------------------------
# Sheet indices (tab order): 1 Sit_Ups, 2 Push_Ups, 3 Reverse_Leg_Lift,
# 4 Shoulder_Press, 5 Squats, 6 Standing_Lunges, 7 Dumbbell_Side_Bend,
# 8 Dumbbell_Curls, 9 Exercise Table

# Reverse_Leg_Lift: new row 4
$wsRevLegLift = $wb.Worksheets.Item(3)
$wsRevLegLift.Range("A4").Value = "25/04/2019"
$wsRevLegLift.Range("B4").Value = 5

# Squats: new row 4
$wsSquats = $wb.Worksheets.Item(5)
$wsSquats.Range("A4").Value = "25/04/2019"
$wsSquats.Range("B4").Value = 15

# Dumbbell_Side_Bend: new row 3
$wsSideBend = $wb.Worksheets.Item(7)
$wsSideBend.Range("A3").Value = "25/04/2019"
$wsSideBend.Range("B3").Value = 15

# Dumbbell_Curls: fill in existing (previously blank-dated) row 4
$wsCurls = $wb.Worksheets.Item(8)
$wsCurls.Range("A4").Value = "25/04/2019"
$wsCurls.Range("B4").Value = 19
$null = $wsCurls.Range("B5").Select()

# --- Make Reverse_Leg_Lift the active / selected sheet & tab --------------
$wsRevLegLift.Activate()

# --- Update the "Exercise Table" summary sheet -----------------------------
$wsTable = $wb.Worksheets.Item(9)
$wsTable.Range("C3").Value = 15
$wsTable.Range("D3").Value = 14
$wsTable.Range("C4").Value = 5
$wsTable.Range("D5").Value = 15
$wsTable.Range("D6").Value = 19
